$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "test case description" column (M2:M40) which contains the
# reference to the URS spec document: bump the version from 1.65 (docx)
# to 1.64 (DOCX, uppercase extension).
$oldText = "製作依據之需求規格書與版本：PJ201800012_URS_5管理性作業_V1.65.docx"
$newText = "製作依據之需求規格書與版本：PJ201800012_URS_5管理性作業_V1.64.DOCX"

for ($r = 2; $r -le 40; $r++) {
    $cell = $ws.Cells.Item($r, 13)  # column M
    if ($cell.Value() -eq $oldText) {
        $cell.Value = $newText
    }
}

# Update the selected cell on the sheet from R2 to A2.
$ws.Range("A2").Select()
